$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 13.75353294838964
$ws.Range("D2").Value = 9.314198999495279
$ws.Range("E2").Value = 14.48201513638037
$ws.Range("F2").Value = 34.62154931655703
$ws.Range("G2").Value = 36.98690962499956
$ws.Range("H2").Value = 16.39164745687302
$ws.Range("J2").Value = 10.64672197330557
$ws.Range("K2").Value = 16.07329673598626
$ws.Range("L2").Value = 9.820060085541362
$ws.Range("N2").Value = 17.6996484220083
$ws.Range("O2").Value = 26.04231660848015
# Row 3
$ws.Range("C3").Value = 13.70609396356384
$ws.Range("D3").Value = 9.285911085735648
$ws.Range("E3").Value = 14.47874736513846
$ws.Range("F3").Value = 34.68471398523327
$ws.Range("G3").Value = 37.05948513017547
$ws.Range("H3").Value = 16.44258800322698
$ws.Range("J3").Value = 10.66660163559762
$ws.Range("K3").Value = 15.66940166720202
$ws.Range("L3").Value = 9.832846465016145
$ws.Range("N3").Value = 17.72399907107707
$ws.Range("O3").Value = 26.12115368454332
# Row 4
$ws.Range("C4").Value = 13.67969289713629
$ws.Range("D4").Value = 9.269768778961424
$ws.Range("E4").Value = 14.47893896849952
$ws.Range("F4").Value = 34.73151786678594
$ws.Range("G4").Value = 37.11501606054799
$ws.Range("H4").Value = 16.47657353105566
$ws.Range("J4").Value = 10.67996885251303
$ws.Range("K4").Value = 15.41732746090355
$ws.Range("L4").Value = 9.841540304724107
$ws.Range("N4").Value = 17.74071860993352
$ws.Range("O4").Value = 26.17511972698597
# Row 5
$ws.Range("C5").Value = 13.66962821145371
$ws.Range("D5").Value = 9.263503590081774
$ws.Range("E5").Value = 14.4795714443479
$ws.Range("F5").Value = 34.75260383523527
$ws.Range("G5").Value = 37.14039408470633
$ws.Range("H5").Value = 16.49110345597783
$ws.Range("J5").Value = 10.68570840027189
$ws.Range("K5").Value = 15.31371889923233
$ws.Range("L5").Value = 9.845295521917784
$ws.Range("N5").Value = 17.74797742032359
$ws.Range("O5").Value = 26.19850641527728
# Row 6
$ws.Range("C6").Value = 13.66799911981204
$ws.Range("D6").Value = 9.262482286703289
$ws.Range("E6").Value = 14.47970999560207
$ws.Range("F6").Value = 34.75622660211795
$ws.Range("G6").Value = 37.14477379938159
$ws.Range("H6").Value = 16.49355723164244
$ws.Range("J6").Value = 10.68667911203048
$ws.Range("K6").Value = 15.29646548242186
$ws.Range("L6").Value = 9.84593191451415
$ws.Range("N6").Value = 17.74920966981018
$ws.Range("O6").Value = 26.20247392781907
# Row 7
$ws.Range("C7").Value = 13.67955434100485
$ws.Range("D7").Value = 9.269683011635729
$ws.Range("E7").Value = 14.4789452514927
$ws.Range("F7").Value = 34.73179409443223
$ws.Range("G7").Value = 37.11534720155525
$ws.Range("H7").Value = 16.47676673168845
$ws.Range("J7").Value = 10.6800450740974
$ws.Range("K7").Value = 15.41593355776136
$ws.Range("L7").Value = 9.841590088288486
$ws.Range("N7").Value = 17.74081470001491
$ws.Range("O7").Value = 26.17542948341357
# Row 8
$ws.Range("C8").Value = 13.73661528856746
$ws.Range("D8").Value = 9.304193524222212
$ws.Range("E8").Value = 14.48043327889356
$ws.Range("F8").Value = 34.64166170474426
$ws.Range("G8").Value = 37.00965155651483
$ws.Range("H8").Value = 16.40864953629671
$ws.Range("J8").Value = 10.65333572211165
$ws.Range("K8").Value = 15.93496232064462
$ws.Range("L8").Value = 9.824294140226382
$ws.Range("N8").Value = 17.70767805871353
$ws.Range("O8").Value = 26.06834396510539
# Row 9
$ws.Range("C9").Value = 13.86974707313426
$ws.Range("D9").Value = 9.381392705804807
$ws.Range("E9").Value = 14.50071025341341
$ws.Range("F9").Value = 34.52870028442837
$ws.Range("G9").Value = 36.88979761568024
$ws.Range("H9").Value = 16.29657309353901
$ws.Range("J9").Value = 10.61015548962931
$ws.Range("K9").Value = 16.91464497343992
$ws.Range("L9").Value = 9.797045907932484
$ws.Range("N9").Value = 17.65668974949245
$ws.Range("O9").Value = 25.90259134886737
# Row 10
$ws.Range("C10").Value = 13.97992251601251
$ws.Range("D10").Value = 9.443623145789822
$ws.Range("E10").Value = 14.52607004920644
$ws.Range("F10").Value = 34.48476453269959
$ws.Range("G10").Value = 36.85547295060857
$ws.Range("H10").Value = 16.22735940451775
$ws.Range("J10").Value = 10.58401652863535
$ws.Range("K10").Value = 17.60406166170281
$ws.Range("L10").Value = 9.781066171805897
$ws.Range("N10").Value = 17.62770988892597
$ws.Range("O10").Value = 25.80795108132751
# Row 11
$ws.Range("C11").Value = 14.03259051033723
$ws.Range("D11").Value = 9.473062738657246
$ws.Range("E11").Value = 14.53984851024621
$ws.Range("F11").Value = 34.47328114640749
$ws.Range("G11").Value = 36.85159230781508
$ws.Range("H11").Value = 16.19872752106631
$ws.Range("J11").Value = 10.57333350782563
$ws.Range("K11").Value = 17.909773536985
$ws.Range("L11").Value = 9.774668071472709
$ws.Range("N11").Value = 17.61635728166806
$ws.Range("O11").Value = 25.77082226338806
# Row 12
$ws.Range("C12").Value = 14.05288834318724
$ws.Range("D12").Value = 9.484367249576284
$ws.Range("E12").Value = 14.54538550431957
$ws.Range("F12").Value = 34.47015611625698
$ws.Range("G12").Value = 36.85181325297754
$ws.Range("H12").Value = 16.18829612122461
$ws.Range("J12").Value = 10.56946141285196
$ws.Range("K12").Value = 18.02430255722103
$ws.Range("L12").Value = 9.772370071375077
$ws.Range("N12").Value = 17.61232066324523
$ws.Range("O12").Value = 25.75761674872966
# Row 13
$ws.Range("C13").Value = 14.04850132610494
$ws.Range("D13").Value = 9.481925767045615
$ws.Range("E13").Value = 14.54417885933006
$ws.Range("F13").Value = 34.4707747289935
$ws.Range("G13").Value = 36.85169045036956
$ws.Range("H13").Value = 16.19052442827433
$ws.Range("J13").Value = 10.57028763419448
$ws.Range("K13").Value = 17.99969322505485
$ws.Range("L13").Value = 9.772859442590413
$ws.Range("N13").Value = 17.61317836610866
$ws.Range("O13").Value = 25.76042275574302
# Row 14
$ws.Range("C14").Value = 14.03425341426649
$ws.Range("D14").Value = 9.47398966985352
$ws.Range("E14").Value = 14.54029765920779
$ws.Range("F14").Value = 34.47299952636618
$ws.Range("G14").Value = 36.85157659297032
$ws.Range("H14").Value = 16.1978610861852
$ws.Range("J14").Value = 10.57301147634398
$ws.Range("K14").Value = 17.91922125260855
$ws.Range("L14").Value = 9.774476514616657
$ws.Range("N14").Value = 17.61601993333938
$ws.Range("O14").Value = 25.76971870102532
# Row 15
$ws.Range("C15").Value = 14.02557181480522
$ws.Range("D15").Value = 9.469148758180864
$ws.Range("E15").Value = 14.53796181218774
$ws.Range("F15").Value = 34.47452162130303
$ws.Range("G15").Value = 36.8517270681718
$ws.Range("H15").Value = 16.20240852348531
$ws.Range("J15").Value = 10.57470247127609
$ws.Range("K15").Value = 17.86976586438514
$ws.Range("L15").Value = 9.775483259891079
$ws.Range("N15").Value = 17.61779461776739
$ws.Range("O15").Value = 25.77552407210699
# Row 16
$ws.Range("C16").Value = 13.97653073551938
$ws.Range("D16").Value = 9.441721424578557
$ws.Range("E16").Value = 14.52521445865701
$ws.Range("F16").Value = 34.48568618404326
$ws.Range("G16").Value = 36.85596297722571
$ws.Range("H16").Value = 16.22928803406773
$ws.Range("J16").Value = 10.58473896313695
$ws.Range("K16").Value = 17.5839147786237
$ws.Range("L16").Value = 9.781501794388138
$ws.Range("N16").Value = 17.62848856980497
$ws.Range("O16").Value = 25.81049695714977
# Row 17
$ws.Range("C17").Value = 13.94708940562757
$ws.Range("D17").Value = 9.42518087813847
$ws.Range("E17").Value = 14.5179665592349
$ws.Range("F17").Value = 34.49471390442431
$ws.Range("G17").Value = 36.8615695917825
$ws.Range("H17").Value = 16.24650902611875
$ws.Range("J17").Value = 10.59120511655589
$ws.Range("K17").Value = 17.40645497670082
$ws.Range("L17").Value = 9.785416777767345
$ws.Range("N17").Value = 17.63551718969992
$ws.Range("O17").Value = 25.83347073224536
# Row 18
$ws.Range("C18").Value = 13.93039623003411
$ws.Range("D18").Value = 9.415774000966994
$ws.Range("E18").Value = 14.51400899823486
$ws.Range("F18").Value = 34.50070677439704
$ws.Range("G18").Value = 36.86589871470853
$ws.Range("H18").Value = 16.25668266644603
$ws.Range("J18").Value = 10.59503797429691
$ws.Range("K18").Value = 17.3036459682981
$ws.Range("L18").Value = 9.787750597940772
$ws.Range("N18").Value = 17.63973221093082
$ws.Range("O18").Value = 25.84724204754131
# Row 19
$ws.Range("C19").Value = 13.92478593091132
$ws.Range("D19").Value = 9.412607524163285
$ws.Range("E19").Value = 14.51270540512997
$ws.Range("F19").Value = 34.50287328649019
$ws.Range("G19").Value = 36.86755403384567
$ws.Range("H19").Value = 16.26017339972467
$ws.Range("J19").Value = 10.59635525178769
$ws.Range("K19").Value = 17.26871291204132
$ws.Range("L19").Value = 9.788554890869333
$ws.Range("N19").Value = 17.64118897087136
$ws.Range("O19").Value = 25.85200043215418
# Row 20
$ws.Range("C20").Value = 13.95019866556098
$ws.Range("D20").Value = 9.426930637698188
$ws.Range("E20").Value = 14.5187162692758
$ws.Range("F20").Value = 34.49367004591652
$ws.Range("G20").Value = 36.86085843483288
$ws.Range("H20").Value = 16.24464802125746
$ws.Range("J20").Value = 10.5905050186458
$ws.Range("K20").Value = 17.42542309572656
$ws.Range("L20").Value = 9.784991535033823
$ws.Range("N20").Value = 17.63475115018283
$ws.Range("O20").Value = 25.83096742584219
# Row 21
$ws.Range("C21").Value = 14.03842888004964
$ws.Range("D21").Value = 9.47631650046471
$ws.Range("E21").Value = 14.54142901890847
$ws.Range("F21").Value = 34.47231284237629
$ws.Range("G21").Value = 36.85156413925678
$ws.Range("H21").Value = 16.19569497699811
$ws.Range("J21").Value = 10.57220671638878
$ws.Range("K21").Value = 17.94289216222475
$ws.Range("L21").Value = 9.773998157367942
$ws.Range("N21").Value = 17.61517818244921
$ws.Range("O21").Value = 25.76696504989011
# Row 22
$ws.Range("C22").Value = 14.09814715096859
$ws.Range("D22").Value = 9.5095014767805
$ws.Range("E22").Value = 14.55813319484056
$ws.Range("F22").Value = 34.46548586940996
$ws.Range("G22").Value = 36.85534348813288
$ws.Range("H22").Value = 16.16609643748186
$ws.Range("J22").Value = 10.56125787398446
$ws.Range("K22").Value = 18.2738348966429
$ws.Range("L22").Value = 9.767540691203594
$ws.Range("N22").Value = 17.60391498968713
$ws.Range("O22").Value = 25.7301167303424
# Row 23
$ws.Range("C23").Value = 14.06609077813219
$ws.Range("D23").Value = 9.491709009079578
$ws.Range("E23").Value = 14.54904872475111
$ws.Range("F23").Value = 34.46847699641457
$ws.Range("G23").Value = 36.85242410599066
$ws.Range("H23").Value = 16.18167442044854
$ws.Range("J23").Value = 10.5670091612721
$ws.Range("K23").Value = 18.09789913989503
$ws.Range("L23").Value = 9.770920761394649
$ws.Range("N23").Value = 17.60978676108686
$ws.Range("O23").Value = 25.74932681844622
# Row 24
$ws.Range("C24").Value = 13.94879224230228
$ws.Range("D24").Value = 9.426139251656206
$ws.Range("E24").Value = 14.51837667303565
$ws.Range("F24").Value = 34.49413947399719
$ws.Range("G24").Value = 36.86117650461369
$ws.Range("H24").Value = 16.24548853100423
$ws.Range("J24").Value = 10.59082117348755
$ws.Range("K24").Value = 17.41685004649153
$ws.Range("L24").Value = 9.785183528567206
$ws.Range("N24").Value = 17.63509693403697
$ws.Range("O24").Value = 25.83209741592994
# Row 25
$ws.Range("C25").Value = 13.83151528693423
$ws.Range("D25").Value = 9.359518009492508
$ws.Range("E25").Value = 14.49337789411797
$ws.Range("F25").Value = 34.55241026759843
$ws.Range("G25").Value = 36.91281214293927
$ws.Range("H25").Value = 16.3245886624497
$ws.Range("J25").Value = 10.6208544556122
$ws.Range("K25").Value = 16.65446485805984
$ws.Range("L25").Value = 9.803705977439870
$ws.Range("N25").Value = 17.66899045525572
$ws.Range("O25").Value = 25.94267753962627
